$wb = $excel.ActiveWorkbook

# --- Sheet "Главные" (index 2) ---
$ws2 = $wb.Worksheets.Item(2)

# Row 2
$ws2.Range("C2").Value = 28
$ws2.Range("D2").Value = 611
$ws2.Range("E2").Value = 253
$ws2.Range("F2").Value = 358
$ws2.Range("G2").Value = 21.82
$ws2.Range("H2").Value = 9.04
$ws2.Range("I2").Value = 12.79
$ws2.Range("J2").Value = 109
$ws2.Range("K2").Value = 134

# Row 5
$ws2.Range("C5").Value = 28
$ws2.Range("D5").Value = 459
$ws2.Range("E5").Value = 234
$ws2.Range("F5").Value = 225
$ws2.Range("G5").Value = 16.39
$ws2.Range("H5").Value = 8.36
$ws2.Range("I5").Value = 8.04
$ws2.Range("J5").Value = 112
$ws2.Range("K5").Value = 105
$ws2.Range("Y5").Value = 4

# Row 8
$ws2.Range("C8").Value = 25
$ws2.Range("D8").Value = 448
$ws2.Range("E8").Value = 225
$ws2.Range("F8").Value = 223
$ws2.Range("G8").Value = 17.92
$ws2.Range("H8").Value = 9
$ws2.Range("I8").Value = 8.92
$ws2.Range("J8").Value = 105
$ws2.Range("K8").Value = 104

# Row 10
$ws2.Range("C10").Value = 18
$ws2.Range("D10").Value = 314
$ws2.Range("E10").Value = 158
$ws2.Range("F10").Value = 156
$ws2.Range("G10").Value = 17.44
$ws2.Range("H10").Value = 8.78
$ws2.Range("I10").Value = 8.67
$ws2.Range("J10").Value = 79
$ws2.Range("K10").Value = 68

# Row 16
$ws2.Range("C16").Value = 27
$ws2.Range("D16").Value = 501
$ws2.Range("E16").Value = 248
$ws2.Range("F16").Value = 253
$ws2.Range("G16").Value = 18.56
$ws2.Range("H16").Value = 9.19
$ws2.Range("I16").Value = 9.37
$ws2.Range("J16").Value = 94
$ws2.Range("K16").Value = 94

# Row 17
$ws2.Range("C17").Value = 17
$ws2.Range("D17").Value = 276
$ws2.Range("E17").Value = 101
$ws2.Range("F17").Value = 175
$ws2.Range("G17").Value = 16.24
$ws2.Range("H17").Value = 5.94
$ws2.Range("I17").Value = 10.29
$ws2.Range("J17").Value = 48
$ws2.Range("K17").Value = 70

# Row 20
$ws2.Range("C20").Value = 26
$ws2.Range("D20").Value = 438
$ws2.Range("E20").Value = 188
$ws2.Range("F20").Value = 250
$ws2.Range("G20").Value = 16.85
$ws2.Range("H20").Value = 7.23
$ws2.Range("I20").Value = 9.62
$ws2.Range("J20").Value = 89
$ws2.Range("K20").Value = 95

# Update as_of_utc timestamp for all data rows (2-26)
for ($r = 2; $r -le 26; $r++) {
    $ws2.Range("AA$r").Value = "2025-11-23 01:57:05"
}

# --- Sheet "Линейные" (index 3) ---
$ws3 = $wb.Worksheets.Item(3)

# Row 3
$ws3.Range("C3").Value = 27
$ws3.Range("D3").Value = 395
$ws3.Range("E3").Value = 203
$ws3.Range("F3").Value = 192
$ws3.Range("G3").Value = 14.63
$ws3.Range("H3").Value = 7.52
$ws3.Range("I3").Value = 7.11
$ws3.Range("J3").Value = 99
$ws3.Range("K3").Value = 81

# Row 8
$ws3.Range("C8").Value = 26
$ws3.Range("D8").Value = 399
$ws3.Range("E8").Value = 150
$ws3.Range("F8").Value = 249
$ws3.Range("G8").Value = 15.35
$ws3.Range("H8").Value = 5.77
$ws3.Range("I8").Value = 9.58
$ws3.Range("J8").Value = 70
$ws3.Range("K8").Value = 97

# Row 9
$ws3.Range("C9").Value = 26
$ws3.Range("D9").Value = 474
$ws3.Range("E9").Value = 203
$ws3.Range("F9").Value = 271
$ws3.Range("G9").Value = 18.23
$ws3.Range("H9").Value = 7.81
$ws3.Range("I9").Value = 10.42
$ws3.Range("J9").Value = 89
$ws3.Range("K9").Value = 113

# Row 14
$ws3.Range("C14").Value = 26
$ws3.Range("D14").Value = 428
$ws3.Range("E14").Value = 218
$ws3.Range("F14").Value = 210
$ws3.Range("G14").Value = 16.46
$ws3.Range("H14").Value = 8.38
$ws3.Range("J14").Value = 109
$ws3.Range("K14").Value = 100

# Row 26
$ws3.Range("C26").Value = 25
$ws3.Range("D26").Value = 503
$ws3.Range("E26").Value = 213
$ws3.Range("F26").Value = 290
$ws3.Range("G26").Value = 20.12
$ws3.Range("H26").Value = 8.52
$ws3.Range("I26").Value = 11.6
$ws3.Range("J26").Value = 84
$ws3.Range("K26").Value = 90
$ws3.Range("Y26").Value = 9

# Update as_of_utc timestamp for all data rows (2-26)
for ($r = 2; $r -le 26; $r++) {
    $ws3.Range("AA$r").Value = "2025-11-23 01:57:05"
}
